# Update cryptos list with newly scraped coinranking.com figures
# (mirrors automated GitHub Actions refresh of cryptos.xlsx)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new "Price" text would otherwise be auto-parsed by Excel as a
# number (losing the exact displayed text, e.g. trailing zeros). Force these
# to Text format before writing the new value so the literal string sticks.
$forceTextCells = @(
  "D4","D5","D6","D8","D9","D10","D13","D14","D15","D16","D19","D20","D22",
  "D25","D27","D28","D29","D32","D33","D35","D37","D38","D39","D40","D41",
  "D42","D44","D45","D46","D47","D48","D49"
)
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Row 2 ---------------------------------------------------------------
$ws.Range("D2").Value = "37.348.26"

# --- Row 3 ---------------------------------------------------------------
$ws.Range("D3").Value = "2.053.98"
$ws.Range("E3").Value = "  -1.44%  "

# --- Row 4 ---------------------------------------------------------------
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.18%  "

# --- Row 5 ---------------------------------------------------------------
$ws.Range("D5").Value = "231.08"
$ws.Range("E5").Value = "  -1.11%  "

# --- Row 6 ---------------------------------------------------------------
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  -0.72%  "

# --- Row 7 ---------------------------------------------------------------
$ws.Range("E7").Value = "  -0.01%  "

# --- Row 8 ---------------------------------------------------------------
$ws.Range("D8").Value = "57.11"
$ws.Range("E8").Value = "  -3.82%  "

# --- Row 9 ---------------------------------------------------------------
$ws.Range("D9").Value = "0.384"
$ws.Range("E9").Value = "  -2.93%  "

# --- Row 10 --------------------------------------------------------------
$ws.Range("D10").Value = "0.0771"
$ws.Range("E10").Value = "  -2.47%  "

# --- Row 11 --------------------------------------------------------------
$ws.Range("E11").Value = "  +1.22%  "

# --- Row 12 --------------------------------------------------------------
$ws.Range("D12").Value = "2.357.28"
$ws.Range("E12").Value = "  -1.41%  "

# --- Row 13 --------------------------------------------------------------
$ws.Range("D13").Value = "14.65"
$ws.Range("E13").Value = "  -1.17%  "

# --- Row 14 --------------------------------------------------------------
$ws.Range("D14").Value = "20.64"
$ws.Range("E14").Value = "  -2.98%  "

# --- Row 15 --------------------------------------------------------------
$ws.Range("D15").Value = "0.758"
$ws.Range("E15").Value = "  -2.46%  "

# --- Row 16 --------------------------------------------------------------
$ws.Range("D16").Value = "5.27"
$ws.Range("E16").Value = "  -1.72%  "

# --- Row 17 --------------------------------------------------------------
$ws.Range("D17").Value = "2.049.20"
$ws.Range("E17").Value = "  -0.23%  "

# --- Row 18 --------------------------------------------------------------
$ws.Range("D18").Value = "37.303.02"
$ws.Range("E18").Value = "  -1.26%  "

# --- Row 19 --------------------------------------------------------------
$ws.Range("D19").Value = "6.08"
$ws.Range("E19").Value = "  -1.64%  "

# --- Row 20 --------------------------------------------------------------
$ws.Range("D20").Value = "69.73"
$ws.Range("E20").Value = "  -2.53%  "

# --- Row 21 --------------------------------------------------------------
$ws.Range("D21").Value = "0.0₃0822"
$ws.Range("E21").Value = "  -3.61%  "

# --- Row 22 --------------------------------------------------------------
$ws.Range("D22").Value = "226.25"
$ws.Range("E22").Value = "  -0.98%  "

# --- Row 23 --------------------------------------------------------------
$ws.Range("E23").Value = "  +0.07%  "

# --- Row 24 --------------------------------------------------------------
$ws.Range("E24").Value = "  +0.01%  "

# --- Row 25 --------------------------------------------------------------
$ws.Range("D25").Value = "2.32"
$ws.Range("E25").Value = "  -4.04%  "

# --- Row 26 --------------------------------------------------------------
$ws.Range("E26").Value = "  +6.37%  "

# --- Row 27 --------------------------------------------------------------
$ws.Range("D27").Value = "169.83"
$ws.Range("E27").Value = "  -1.28%  "

# --- Row 28 --------------------------------------------------------------
$ws.Range("D28").Value = "0.130"
$ws.Range("E28").Value = "  -5.76%  "

# --- Row 29 --------------------------------------------------------------
$ws.Range("D29").Value = "19.16"
$ws.Range("E29").Value = "  -1.79%  "

# --- Row 30 --------------------------------------------------------------
$ws.Range("E30").Value = "  -5.52%  "

# --- Row 31 --------------------------------------------------------------
$ws.Range("E31").Value = "  -0.01%  "

# --- Row 32 --------------------------------------------------------------
$ws.Range("D32").Value = "4.53"
$ws.Range("E32").Value = "  -4.40%  "

# --- Row 33 --------------------------------------------------------------
$ws.Range("D33").Value = "0.0623"
$ws.Range("E33").Value = "  -1.64%  "

# --- Row 34 --------------------------------------------------------------
$ws.Range("E34").Value = "  -4.44%  "

# --- Row 35 --------------------------------------------------------------
$ws.Range("D35").Value = "2.48"
$ws.Range("E35").Value = "  -1.78%  "

# --- Row 36 --------------------------------------------------------------
$ws.Range("E36").Value = "  -0.01%  "

# --- Row 37 --------------------------------------------------------------
$ws.Range("D37").Value = "3.28"
$ws.Range("E37").Value = "  -4.98%  "

# --- Row 38 --------------------------------------------------------------
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.08%  "

# --- Row 39 --------------------------------------------------------------
$ws.Range("D39").Value = "5.31"
$ws.Range("E39").Value = "  -2.19%  "

# --- Row 40 --------------------------------------------------------------
$ws.Range("D40").Value = "0.0226"
$ws.Range("E40").Value = "  +3.49%  "

# --- Row 41 --------------------------------------------------------------
$ws.Range("D41").Value = "98.37"
$ws.Range("E41").Value = "  -0.97%  "

# --- Row 42 (Maker -> Cronos) ---------------------------------------------
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "0.0953"
$ws.Range("E42").Value = "  -3.38%  "

# --- Row 43 (HuobiToken -> Maker) -----------------------------------------
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.483.26"
$ws.Range("E43").Value = "  +2.54%  "

# --- Row 44 (Cronos -> HuobiToken) -----------------------------------------
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "2.90"
$ws.Range("E44").Value = "  +0.48%  "

# --- Row 45 --------------------------------------------------------------
$ws.Range("D45").Value = "1.18"
$ws.Range("E45").Value = "  +1.52%  "

# --- Row 46 --------------------------------------------------------------
$ws.Range("D46").Value = "16.54"
$ws.Range("E46").Value = "  -0.76%  "

# --- Row 47 (ARBITRUM -> FTXToken) ----------------------------------------
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").Value = "4.00"
$ws.Range("E47").Value = "  -4.71%  "

# --- Row 48 (FTXToken -> ARBITRUM) ----------------------------------------
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "1.03"
$ws.Range("E48").Value = "  -3.10%  "

# --- Row 49 --------------------------------------------------------------
$ws.Range("D49").Value = "7.25"
$ws.Range("E49").Value = "  -2.19%  "

# --- Row 50 --------------------------------------------------------------
$ws.Range("E50").Value = "  -2.13%  "

# --- Row 51 --------------------------------------------------------------
$ws.Range("D51").Value = "2.242.87"
$ws.Range("E51").Value = "  -1.44%  "
